$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 96
$ws.Cells.Item($row, 1).Value = "04-10-2021"
$ws.Cells.Item($row, 2).Value = 50000
$ws.Cells.Item($row, 3).Value = 45000
$ws.Cells.Item($row, 4).Value = 45000
$ws.Cells.Item($row, 5).Value = 40000
$ws.Cells.Item($row, 6).Value = 5000
$ws.Cells.Item($row, 7).Value = 2.5
